# Update the "取得日時" (retrieved datetime) timestamps on the first sheet
# (ランサーズ) for rows 2-7 from 2025-09-21 12:32:58 to 2025-09-21 12:40:50.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 1).Value = "2025-09-21 12:40:50"
}
